$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns widths
$ws.Columns.Item(4).ColumnWidth = 23.88671875
$ws.Columns.Item(5).ColumnWidth = 18.77734375

# Enter values in the same sequence a human editor would have typed them,
# so that new shared-string entries are created in that order.
$ws.Range("D1").Value = "After spell check"
$ws.Range("D2").Value = "128/149 (85.9%)"

$ws.Range("A3").Value = "Arial"
$ws.Range("B3").Value = "211/243"
$ws.Range("C3").Value = 0.868

$ws.Range("E2").Value = "Neural Network"
$ws.Range("E1").Value = "Classification Method"
$ws.Range("E3").Value = "Least Distance"
$ws.Range("D3").Value = "230/243 (94.7%)"
$ws.Range("F3").Value = "*** 20x20 raw img input w/ sets 2 and 3"

$ws.Range("A4").Value = "Arial"
$ws.Range("B4").Value = "224/243"
$ws.Range("C4").Value = 0.922
$ws.Range("D4").Value = "224/243 (92.2%)"
$ws.Range("E4").Value = "Least Distance"

$ws.Range("A5").Value = "Times New Roman"
$ws.Range("B5").Value = "83/117"
$ws.Range("C5").Value = 0.709
$ws.Range("D5").Value = "N/A"
$ws.Range("E5").Value = "Least Distance"

# Number formats (percentage style, matching existing C2 cell) - copy the
# single formatted cell so no new cellXfs entries are minted.
$ws.Range("C2").Copy()
$ws.Range("C3:C5").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header formatting for D1:E1 to match A1 style (single-cell copy avoids
# tiling past the destination and avoids minting new cellXfs entries)
$ws.Range("A1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B8").Select()
